# Updated symbol list on Sat Feb 11 21:57:15 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns of the crypto tracker sheet.
# Values are stored as literal text (e.g. "309.90", "0.88%") in the source
# workbook, so each cell is forced to Text format before the write -- otherwise
# Excel would silently reinterpret "0.88%" as the number 0.0088, or drop the
# trailing zero in "309.90".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '309.90'
    'E2' = '0.88%'
    'D3' = '40.82'
    'E3' = '0.89%'
    'D4' = '5.114'
    'E4' = '1.21%'
    'D5' = '0.07667'
    'E5' = '1.04%'
    'D6' = '1.613'
    'E6' = '0.74%'
    'D7' = '0.9145'
    'E7' = '0.83%'
    'E8' = '0.67%'
    'D9' = '0.1261'
    'E9' = '23.57%'
    'D10' = '0.1804'
    'E10' = '3.06%'
    'D11' = '0.09133'
    'E11' = '0.54%'
    'D12' = '0.04291'
    'E12' = '1.98%'
    'D13' = '0.1048'
    'E13' = '-0.37%'
    'D14' = '0.001229'
    'E14' = '-2.16%'
    'D15' = '0.005774'
    'E15' = '-0.76%'
    'D16' = '3.351'
    'E16' = '-0.04%'
    'D17' = '4.276'
    'E17' = '0.29%'
    'E18' = '1.54%'
    'D19' = '6.895'
    'E19' = '1.93%'
    'D20' = '0.1390'
    'E20' = '1.74%'
    'D22' = '0.04029'
    'E22' = '-3.75%'
    'D23' = '0.001267'
    'E23' = '3.21%'
    'D24' = '0.004141'
    'E24' = '2.38%'
    'D25' = '0.0001268'
    'E25' = '-2.79%'
    'E26' = '24.47%'
    'D38' = '0.02435'
    'E38' = '2.01%'
    'D39' = '0.05262'
    'E39' = '1.86%'
    'D40' = '0.007826'
    'E40' = '0.56%'
    'D41' = '0.1310'
    'E41' = '0.94%'
    'D42' = '0.006787'
    'E42' = '-4.02%'
    'D43' = '0.001839'
    'E43' = '-4.47%'
    'D44' = '0.008196'
    'E44' = '-3.79%'
    'D45' = '0.3348'
    'E45' = '0.01%'
    'D46' = '0.00006844'
    'E46' = '7.43%'
    'D47' = '0.00000000748'
    'E47' = '-0.61%'
    'D48' = '0.3193'
    'E48' = '4,398.59%'
    'E49' = '-31.91%'
    'D50' = '0.00002095'
    'E50' = '-0.61%'
    'D51' = '0.0001995'
    'E51' = '-0.61%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
